$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the source schema on row 2 (was a leftover/typo value)
$ws.Range("C2").Value = "USAG_LOCATION_SYNC"

# Correct the target table reference on row 9 (should be T1, matching row 7/8)
$ws.Range("D9").Value = "T1"
